# Commit: "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The account-statement worker table (rows 16-62) is updated:
#  - Rows 16/17 (MARCO FIDEL SOSA VENECIA <-> ROGELIO HERNANDEZ PEREZ) swap
#    places, so their Doc# / Name / Valor Mora / Salario Basico trade rows.
#  - Row 18 (JANE PAOLA CARVAJAL ORTIZ) is untouched.
#  - Rows 19-62 (JOSE GREGORIO RICO OROZCO, one row per overdue period) have
#    their "Periodo Mora" list reversed (was ascending 1703..2012, now
#    descending 2012..1703), the matching "Valor Mora" value travels with its
#    original period (same reversal), and "Salario Basico" is unified to the
#    new value 737717 for every one of those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 16 / 17: swap the two workers' identity + amounts ---------------
$ws.Range("C16").Value = "73079608"
$ws.Range("D16").Value = "MARCO FIDEL SOSA VENECIA"
$ws.Range("E16").Value = "1703"
$ws.Range("F16").Value = 60000
$ws.Range("G16").Value = 1500000

$ws.Range("C17").Value = "73352479"
$ws.Range("D17").Value = "ROGELIO HERNANDEZ PEREZ"
$ws.Range("E17").Value = "1703"
$ws.Range("F17").Value = 29509
$ws.Range("G17").Value = 737717

# --- Row 18: unchanged (kept explicit for clarity) -------------------------
$ws.Range("E18").Value = "1703"
$ws.Range("F18").Value = 60000
$ws.Range("G18").Value = 1500000

# --- Rows 19-62: JOSE GREGORIO RICO OROZCO, periods reversed --------------
# Original (row19..row62) period order, ascending:
$periods = @("1703","1706","1707","1708","1709","1710","1711","1712","1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812","1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012")

# Original (row19..row62) "Valor Mora" values, matched index-for-index with $periods above:
$valorMora = @(27578,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,20656)

if ($periods.Length -ne 44 -or $valorMora.Length -ne 44) {
    throw "unexpected array length: periods=$($periods.Length) valorMora=$($valorMora.Length)"
}

$newSalarioBasico = 737717
$startRow = 19
$count = $periods.Length

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $srcIndex = $count - 1 - $i
    $ws.Range("E$row").Value = $periods[$srcIndex]
    $ws.Range("F$row").Value = $valorMora[$srcIndex]
    $ws.Range("G$row").Value = $newSalarioBasico
}
